# Data source documentation update
# Adds three new plain-text URL notes to the "Data Source Spreadsheet":
#   K3 -> NCES school district finance SAS7BDAT download link
#   I6 -> COVID-19 School Data Hub white paper PDF
#   I9 -> Senator Lee social capital data file link
# Also widens column J (10) to fit the new note text and tightens row 3's
# height now that its content fits more compactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = "https://nces.ed.gov/ccd/Data/zip/sdf22_1a_sas7bdat.zip"
$ws.Range("I6").Value = "https://assets.ctfassets.net/9fbw4onh0qc1/51GumHMBAalCkUBBHREfxO/0f4e84f44f90451b5558e7989e0d342b/CSDH_Introductory_White_Paper.pdf"
$ws.Range("I9").Value = "https://www.lee.senate.gov/services/files/DA64FDB7-3B2E-40D4-B9E3-07001B81EC31"

# New column J width to comfortably show the added note text.
$ws.Columns.Item(10).ColumnWidth = 18

# Row 3 shrank after review (K3 note fits without extra wrapped lines).
$ws.Rows.Item(3).RowHeight = 105

# Restore the view: less zoomed in, scrolled down so row 8 pins under the
# frozen header row, with the active cell resting on the newly-added J10.
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("A8").Select() | Out-Null
$ws.Range("J10").Select() | Out-Null
